# Updated cryptos list values (Price and Volume(1h) columns) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.183.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.71%  "

$ws.Range("D3").Value = "'2.633.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.01%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'606.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.18%  "

$ws.Range("D6").Value = "'181.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.99%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +1.14%  "

$ws.Range("D9").Value = "'2.631.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.92%  "

$ws.Range("E10").Value = "  +13.69%  "

$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").Value = "'0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.02%  "

$ws.Range("D13").Value = "'5.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D14").Value = "'0.0000189"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.72%  "

$ws.Range("D15").Value = "'3.051.42"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'26.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.08%  "

$ws.Range("D17").Value = "'71.157.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.85%  "

$ws.Range("D18").Value = "'2.618.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.46%  "

$ws.Range("D19").Value = "'382.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.20%  "

$ws.Range("D20").Value = "'7.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.77%  "

$ws.Range("D21").Value = "'11.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.62%  "

$ws.Range("E22").Value = "  -1.81%  "

$ws.Range("D23").Value = "'4.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.00%  "

$ws.Range("D24").Value = "'72.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.84%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("E26").Value = "  +11.36%  "

$ws.Range("D27").Value = "'9.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.68%  "

$ws.Range("D28").Value = "'2.766.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.15%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "'0.0₃0965"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.37%  "

$ws.Range("D31").Value = "'541.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.54%  "

$ws.Range("D32").Value = "'8.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.00%  "

$ws.Range("D33").Value = "'1.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.11%  "

$ws.Range("E34").Value = "  +3.27%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").Value = "'165.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.04%  "

$ws.Range("D37").Value = "'0.118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.59%  "

$ws.Range("D38").Value = "'19.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.18%  "

$ws.Range("E39").Value = "  +7.39%  "

$ws.Range("E40").Value = "  +1.84%  "

$ws.Range("E41").Value = "  +4.95%  "

$ws.Range("D42").Value = "'2.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.21%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").Value = "'5.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.98%  "

$ws.Range("E45").Value = "  +1.18%  "

$ws.Range("D46").Value = "'39.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.50%  "

$ws.Range("D47").Value = "'154.21"
$ws.Range("D47").Style = "Normal"

$ws.Range("E48").Value = "  +1.76%  "

$ws.Range("E49").Value = "  +5.14%  "

$ws.Range("E50").Value = "  +2.23%  "

$ws.Range("D51").Value = "'0.0₆0264"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.30%  "
